# "added autodetect row format and combine columns"
#
# Adds 5 new "indent" rows (indent_1..indent_5) below the existing table,
# each progressively smaller in font size and with a progressively darker
# accent-4 themed fill (the standard Excel "Lighter 80/60/40%, Darker 25/50%"
# ramp), then moves the selection to E9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 11: indent_1 -------------------------------------------------
$ws.Range("A11").Value = "indent_1"
$ws.Range("A11").Font.Size = 14
$ws.Range("A11").Interior.Color = 14942438   # accent4, Lighter 80% (#E6E0EC)
$ws.Range("B11").Value = 15
$ws.Rows.Item(11).RowHeight = 18

# --- row 13: indent_2 -------------------------------------------------
$ws.Range("A13").Value = "indent_2"
$ws.Range("A13").Interior.Color = 14333900   # accent4, Lighter 60% (#CCC1DA)
$ws.Range("B13").Value = 13

# --- row 15: indent_3 -------------------------------------------------
$ws.Range("A15").Value = "indent_3"
$ws.Range("A15").Font.Size = 10
$ws.Range("A15").Interior.Color = 13085363   # accent4, Lighter 40% (#B3A2C7)
$ws.Range("B15").Value = 11

# --- row 17: indent_4 -------------------------------------------------
$ws.Range("A17").Value = "indent_4"
$ws.Range("A17").Font.Size = 8
$ws.Range("A17").Font.ThemeColor = 2         # background1 / white text
$ws.Range("A17").Interior.Color = 8095616    # accent4, Darker 25% (#604A7B)
$ws.Range("B17").Value = 9

# --- row 19: indent_5 -------------------------------------------------
$ws.Range("A19").Value = "indent_5"
$ws.Range("A19").Font.Size = 6
$ws.Range("A19").Font.ThemeColor = 2         # background1 / white text
$ws.Range("A19").Interior.Color = 4210002    # accent4, Darker 50% (#403152)
$ws.Range("B19").Value = 7

$ws.Range("E9").Select()
